$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 0.1672597864768683
$ws.Range("C2").Value = 0.608540925266904
$ws.Range("J2").Value = 0.003558718861209964
$ws.Range("P2").Value = 0.1209964412811388
$ws.Range("S2").Value = 0.099644128113879
$ws.Range("C3").Value = 0.02259887005649718
$ws.Range("P3").Value = 0.8305084745762712
$ws.Range("S3").Value = 0.1468926553672316
$ws.Range("J4").Value = 0.01724137931034483
$ws.Range("P4").Value = 0.7413793103448276
$ws.Range("S4").Value = 0.2413793103448276
$ws.Range("B6").Value = 0.06220095693779904
$ws.Range("D6").Value = 0.01435406698564593
$ws.Range("F6").Value = 0.05263157894736842
$ws.Range("J6").Value = 0.2488038277511962
$ws.Range("O6").Value = 0.02392344497607655
$ws.Range("Q6").Value = 0.1674641148325359
$ws.Range("R6").Value = 0.07655502392344497
$ws.Range("S6").Value = 0.354066985645933
$ws.Range("B7").Value = 0.09202453987730061
$ws.Range("D7").Value = 0.049079754601227
$ws.Range("F7").Value = 0.0736196319018405
$ws.Range("J7").Value = 0.147239263803681
$ws.Range("Q7").Value = 0.2024539877300613
$ws.Range("R7").Value = 0.04294478527607362
$ws.Range("S7").Value = 0.392638036809816
$ws.Range("B8").Value = 0.08333333333333333
$ws.Range("D8").Value = 0.009803921568627451
$ws.Range("F8").Value = 0.05392156862745098
$ws.Range("J8").Value = 0.1911764705882353
$ws.Range("O8").Value = 0.02450980392156863
$ws.Range("Q8").Value = 0.2156862745098039
$ws.Range("R8").Value = 0.04656862745098039
$ws.Range("S8").Value = 0.375
$ws.Range("B9").Value = 0.09051724137931035
$ws.Range("D9").Value = 0.01724137931034483
$ws.Range("F9").Value = 0.04741379310344827
$ws.Range("J9").Value = 0.1379310344827586
$ws.Range("O9").Value = 0.02155172413793104
$ws.Range("Q9").Value = 0.1724137931034483
$ws.Range("R9").Value = 0.08620689655172414
$ws.Range("S9").Value = 0.4267241379310345
$ws.Range("B10").Value = 0.1089324618736384
$ws.Range("D10").Value = 0.03050108932461874
$ws.Range("F10").Value = 0.074800290486565
$ws.Range("J10").Value = 0.1263616557734205
$ws.Range("O10").Value = 0.01742919389978214
$ws.Range("Q10").Value = 0.2389251997095134
$ws.Range("R10").Value = 0.06100217864923747
$ws.Range("S10").Value = 0.3420479302832244
$ws.Range("G11").Value = 0.1048387096774194
$ws.Range("J11").Value = 0.125
$ws.Range("K11").Value = 0.1774193548387097
$ws.Range("L11").Value = 0.592741935483871
$ws.Range("G12").Value = 0.7908496732026143
$ws.Range("J12").Value = 0.1568627450980392
$ws.Range("K12").Value = 0.006535947712418301
$ws.Range("L12").Value = 0.03267973856209151
$ws.Range("S12").Value = 0.0130718954248366
$ws.Range("F15").Value = 0.00425531914893617
$ws.Range("H15").Value = 0.1404255319148936
$ws.Range("I15").Value = 0.07234042553191489
$ws.Range("J15").Value = 0.4212765957446808
$ws.Range("K15").Value = 0.05531914893617021
$ws.Range("M15").Value = 0.00425531914893617
$ws.Range("N15").Value = 0.00425531914893617
$ws.Range("O15").Value = 0.03404255319148936
$ws.Range("S15").Value = 0.2638297872340425
$ws.Range("F16").Value = 0.0186046511627907
$ws.Range("H16").Value = 0.1534883720930233
$ws.Range("I16").Value = 0.1162790697674419
$ws.Range("J16").Value = 0.4186046511627907
$ws.Range("K16").Value = 0.06511627906976744
$ws.Range("M16").Value = 0.01395348837209302
$ws.Range("O16").Value = 0.07441860465116279
$ws.Range("S16").Value = 0.1395348837209302
$ws.Range("F17").Value = 0.01149425287356322
$ws.Range("H17").Value = 0.1628352490421456
$ws.Range("I17").Value = 0.1053639846743295
$ws.Range("J17").Value = 0.446360153256705
$ws.Range("K17").Value = 0.0842911877394636
$ws.Range("M17").Value = 0.01724137931034483
$ws.Range("O17").Value = 0.06896551724137931
$ws.Range("S17").Value = 0.103448275862069
$ws.Range("F18").Value = 0.01360544217687075
$ws.Range("H18").Value = 0.108843537414966
$ws.Range("I18").Value = 0.09523809523809523
$ws.Range("J18").Value = 0.3945578231292517
$ws.Range("K18").Value = 0.1020408163265306
$ws.Range("M18").Value = 0.03401360544217687
$ws.Range("O18").Value = 0.1156462585034014
$ws.Range("S18").Value = 0.1360544217687075
$ws.Range("F19").Value = 0.01049233252623083
$ws.Range("H19").Value = 0.1945117029862793
$ws.Range("I19").Value = 0.09685230024213075
$ws.Range("J19").Value = 0.3938660209846651
$ws.Range("K19").Value = 0.09362389023405973
$ws.Range("M19").Value = 0.01372074253430186
$ws.Range("N19").Value = 0.001614205004035512
$ws.Range("O19").Value = 0.06941081517352704
$ws.Range("S19").Value = 0.12590799031477
